# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Spriggan_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled price-refresh diff.

$wb = $excel.ActiveWorkbook

# ========== Sheet: ALC ==========
$ws = $wb.Worksheets.Item("ALC")

# Row 116
$ws.Range("H116").Value = 3661.0386
$ws.Range("I116").Value = 3731.0527
$ws.Range("K116").Value = 3731.0527
$ws.Range("M116").Value = -289.0527000000002

# Row 125
$ws.Range("H125").Value = 3045.9285
$ws.Range("J125").Value = 3166.6667
$ws.Range("L125").Value = 28500.0003
$ws.Range("N125").Value = -33420.0003

# Row 127
$ws.Range("H127").Value = 1708.625
$ws.Range("I127").Value = 1238.4286
$ws.Range("J127").Value = 5000
$ws.Range("K127").Value = 3715.2858
$ws.Range("L127").Value = 15000
$ws.Range("M127").Value = 1244.7142
$ws.Range("N127").Value = -24920

# Row 132
$ws.Range("H132").Value = 1836.6129
$ws.Range("I132").Value = 1836.6129
$ws.Range("K132").Value = 5509.8387
$ws.Range("M132").Value = -2979.8387

# Row 137
$ws.Range("H137").Value = 2323.4
$ws.Range("I137").Value = 2041.3214
$ws.Range("K137").Value = 6123.9642
$ws.Range("M137").Value = -3573.9642

# Row 138
$ws.Range("H138").Value = 2556.2952
$ws.Range("I138").Value = 1455.1
$ws.Range("K138").Value = 4365.299999999999
$ws.Range("M138").Value = 774.7000000000007

# ========== Sheet: ARM ==========
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 2715.5425
$ws.Range("I32").Value = 2758.1914
$ws.Range("K32").Value = 2758.1914
$ws.Range("M32").Value = -2471.1914

# Row 56
$ws.Range("H56").Value = 6000
$ws.Range("I56").Value = 6000
$ws.Range("K56").Value = 6000
$ws.Range("M56").Value = -5258

# Row 74
$ws.Range("H74").Value = 52640900
$ws.Range("I74").Value = 66675876
$ws.Range("J74").Value = 9749.75
$ws.Range("K74").Value = 66675876
$ws.Range("L74").Value = 9749.75
$ws.Range("M74").Value = -66675002
$ws.Range("N74").Value = -11497.75

# Row 77
$ws.Range("H77").Value = 52640900
$ws.Range("I77").Value = 66675876
$ws.Range("J77").Value = 9749.75
$ws.Range("K77").Value = 333379380
$ws.Range("L77").Value = 48748.75
$ws.Range("M77").Value = -333375012
$ws.Range("N77").Value = -57484.75

# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# Row 122
$ws.Range("H122").Value = 2965.2856
$ws.Range("I122").Value = 2855.818
$ws.Range("J122").Value = 3366.6667
$ws.Range("K122").Value = 8567.454000000002
$ws.Range("L122").Value = 10100.0001
$ws.Range("M122").Value = -6117.454000000002
$ws.Range("N122").Value = -15000.0001

# Row 134
$ws.Range("H134").Value = 210000
$ws.Range("J134").Value = 210000
$ws.Range("L134").Value = 210000
$ws.Range("N134").Value = -220140

# Row 137
$ws.Range("H137").Value = 60000
$ws.Range("I137").Value = 60000
$ws.Range("K137").Value = 60000
$ws.Range("M137").Value = -54900

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ========== Sheet: BSM ==========
$ws = $wb.Worksheets.Item("BSM")

# Row 82
$ws.Range("H82").Value = 47499.25

# Row 85
$ws.Range("H85").Value = 47499.25

# ========== Sheet: CRP ==========
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 7996.0527
$ws.Range("J31").Value = 9343.429
$ws.Range("L31").Value = 9343.429
$ws.Range("N31").Value = -9933.429

# Row 34
$ws.Range("H34").Value = 7996.0527
$ws.Range("J34").Value = 9343.429
$ws.Range("L34").Value = 9343.429
$ws.Range("N34").Value = -9747.429

# Row 58
$ws.Range("H58").Value = 15629043
$ws.Range("I58").Value = 20004720
$ws.Range("J58").Value = 1626.2858
$ws.Range("K58").Value = 20004720
$ws.Range("L58").Value = 1626.2858
$ws.Range("M58").Value = -20004517
$ws.Range("N58").Value = -2032.2858

# Row 60
$ws.Range("H60").Value = 29594.95
$ws.Range("I60").Value = 3650
$ws.Range("K60").Value = 3650
$ws.Range("M60").Value = -3139

# Row 122
$ws.Range("H122").Value = 2912.6428
$ws.Range("I122").Value = 2983.5557
$ws.Range("K122").Value = 8950.667099999999
$ws.Range("M122").Value = -6500.667099999999

# Row 134
$ws.Range("H134").Value = 4167965.5
$ws.Range("I134").Value = 4718227
$ws.Range("K134").Value = 14154681
$ws.Range("M134").Value = -14152146

# Row 136
$ws.Range("H136").Value = 15629043
$ws.Range("I136").Value = 20004720
$ws.Range("J136").Value = 1626.2858
$ws.Range("K136").Value = 60014160
$ws.Range("L136").Value = 4878.857400000001
$ws.Range("M136").Value = -60011610
$ws.Range("N136").Value = -9978.857400000001

# ========== Sheet: CUL ==========
$ws = $wb.Worksheets.Item("CUL")

# Row 44
$ws.Range("H44").Value = 309.16666
$ws.Range("I44").Value = 309.16666
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 927.4999799999999
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -529.4999799999999
$ws.Range("N44").ClearContents()

# Row 68
$ws.Range("H68").Value = 4666.3335
$ws.Range("J68").Value = 4999
$ws.Range("L68").Value = 14997
$ws.Range("N68").Value = -16619

# Row 69
$ws.Range("H69").Value = 741.6667
$ws.Range("I69").Value = 741.6667
$ws.Range("K69").Value = 2225.0001
$ws.Range("M69").Value = -1414.0001

# Row 71
$ws.Range("H71").Value = 4666.3335
$ws.Range("J71").Value = 4999
$ws.Range("L71").Value = 44991
$ws.Range("N71").Value = -53103

# Row 72
$ws.Range("H72").Value = 741.6667
$ws.Range("I72").Value = 741.6667
$ws.Range("K72").Value = 6675.0003
$ws.Range("M72").Value = -2619.0003

# Row 119
$ws.Range("H119").Value = 4712.8887
$ws.Range("I119").Value = 2802.125
$ws.Range("J119").Value = 19999
$ws.Range("K119").Value = 8406.375
$ws.Range("L119").Value = 59997
$ws.Range("M119").Value = -3568.375
$ws.Range("N119").Value = -69673

# Row 120
$ws.Range("H120").Value = 20666
$ws.Range("I120").Value = 15999.5
$ws.Range("J120").Value = 29999
$ws.Range("K120").Value = 47998.5
$ws.Range("L120").Value = 89997
$ws.Range("M120").Value = -43160.5
$ws.Range("N120").Value = -99673

# Row 131
$ws.Range("H131").Value = 1600.8235
$ws.Range("I131").Value = 1154.2667
$ws.Range("K131").Value = 3462.800099999999
$ws.Range("M131").Value = 1577.199900000001

# ========== Sheet: GSM ==========
$ws = $wb.Worksheets.Item("GSM")

# Row 5
$ws.Range("H5").Value = 900
$ws.Range("J5").Value = 900
$ws.Range("L5").Value = 900
$ws.Range("N5").Value = -1124

# Row 100
$ws.Range("H100").Value = 129999
$ws.Range("J100").Value = 129999
$ws.Range("L100").Value = 129999
$ws.Range("N100").Value = -132163

# Row 122
$ws.Range("H122").Value = 4436.079
$ws.Range("I122").Value = 2856.2
$ws.Range("K122").Value = 8568.599999999999
$ws.Range("M122").Value = -6118.599999999999

# Row 126
$ws.Range("H126").Value = 8510.333000000001
$ws.Range("I126").Value = 9571.286
$ws.Range("J126").Value = 4797
$ws.Range("K126").Value = 28713.858
$ws.Range("L126").Value = 14391
$ws.Range("M126").Value = -26243.858
$ws.Range("N126").Value = -19331

# Row 132
$ws.Range("H132").Value = 3677534
$ws.Range("I132").Value = 3907317.5
$ws.Range("K132").Value = 11721952.5
$ws.Range("M132").Value = -11719422.5

# Row 136
$ws.Range("H136").Value = 39599.727
$ws.Range("J136").Value = 39599.727
$ws.Range("L136").Value = 118799.181
$ws.Range("N136").Value = -123899.181

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ========== Sheet: LTW ==========
$ws = $wb.Worksheets.Item("LTW")

# Row 61
$ws.Range("H61").Value = 3102
$ws.Range("I61").Value = 3250
$ws.Range("K61").Value = 3250
$ws.Range("M61").Value = -3048

# Row 113
$ws.Range("H113").Value = 3102
$ws.Range("I113").Value = 3250
$ws.Range("K113").Value = 3250
$ws.Range("M113").Value = -1080

# Row 122
$ws.Range("H122").Value = 4797.3076
$ws.Range("I122").Value = 4806.25
$ws.Range("K122").Value = 14418.75
$ws.Range("M122").Value = -11968.75

# Row 136
$ws.Range("H136").Value = 2046.7059
$ws.Range("I136").Value = 1050.2858
$ws.Range("J136").Value = 2744.2
$ws.Range("K136").Value = 3150.8574
$ws.Range("L136").Value = 8232.599999999999
$ws.Range("M136").Value = -600.8574000000003
$ws.Range("N136").Value = -13332.6

# ========== Sheet: WVR ==========
$ws = $wb.Worksheets.Item("WVR")

# Row 74
$ws.Range("H74").Value = 13928.556
$ws.Range("J74").Value = 13060.833
$ws.Range("L74").Value = 13060.833
$ws.Range("N74").Value = -14932.833

# Row 77
$ws.Range("H77").Value = 13928.556
$ws.Range("J77").Value = 13060.833
$ws.Range("L77").Value = 39182.499
$ws.Range("N77").Value = -48542.499

# Row 122
$ws.Range("H122").Value = 1279.7693
$ws.Range("I122").Value = 1279.7693
$ws.Range("K122").Value = 3839.3079
$ws.Range("M122").Value = -1389.3079

# Row 132
$ws.Range("H132").Value = 11631827
$ws.Range("I132").Value = 13890352
$ws.Range("K132").Value = 41671056
$ws.Range("M132").Value = -41668526

